$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("2023-12-14 16:18:23", 0.0006000000000000001),
    @("2023-12-14 16:18:40", 0.0008),
    @("2023-12-14 16:18:58", 0.0012),
    @("2023-12-14 16:19:06", 0.0002)
)

$startRow = 324
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}
